# Updates cryptos list figures (prices in column D, 1h volume % in column E),
# and swaps the TheSandbox / FraxShare row positions (rows 40-41).
# Prices in column D are forced to remain TEXT (leading apostrophe) so that
# number-like strings such as "46.28" are not auto-converted to numeric values,
# matching the original inline-string/text storage of that column.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'28.607.00"
$ws.Range("E2").Value = "  -3.41%  "
$ws.Range("D3").Value = "'1.850.43"
$ws.Range("E3").Value = "  -4.00%  "
$ws.Range("E4").Value = "  -1.20%  "
$ws.Range("D5").Value = "'335.44"
$ws.Range("E5").Value = "  +2.74%  "
$ws.Range("E6").Value = "  -1.07%  "
$ws.Range("D7").Value = "'0.4654"
$ws.Range("E7").Value = "  -3.55%  "
$ws.Range("D8").Value = "'0.3903"
$ws.Range("E8").Value = "  -3.82%  "
$ws.Range("D9").Value = "'46.28"
$ws.Range("E9").Value = "  -2.52%  "
$ws.Range("D10").Value = "'0.07908"
$ws.Range("E10").Value = "  -3.52%  "
$ws.Range("D11").Value = "'0.9773"
$ws.Range("E11").Value = "  -3.05%  "
$ws.Range("D12").Value = "'22.28"
$ws.Range("E12").Value = "  -6.11%  "
$ws.Range("D13").Value = "'1.902.14"
$ws.Range("E13").Value = "  -0.66%  "
$ws.Range("D14").Value = "'5.813"
$ws.Range("E14").Value = "  -4.23%  "
$ws.Range("D15").Value = "'6.964"
$ws.Range("E15").Value = "  -4.43%  "
$ws.Range("D16").Value = "'0.06907"
$ws.Range("E16").Value = "  +0.69%  "
$ws.Range("D17").Value = "'87.78"
$ws.Range("E18").Value = "  -1.14%  "
$ws.Range("D19").Value = "'0.00001002"
$ws.Range("E19").Value = "  -3.50%  "
$ws.Range("E20").Value = "  -3.05%  "
$ws.Range("D21").Value = "'1.001"
$ws.Range("E21").Value = "  -1.04%  "
$ws.Range("D22").Value = "'28.614.04"
$ws.Range("E22").Value = "  -3.35%  "
$ws.Range("D23").Value = "'5.381"
$ws.Range("E23").Value = "  -4.80%  "
$ws.Range("D24").Value = "'11.22"
$ws.Range("E24").Value = "  -5.97%  "
$ws.Range("D25").Value = "'2.151"
$ws.Range("E25").Value = "  -2.14%  "
$ws.Range("D26").Value = "'2.081.30"
$ws.Range("E26").Value = "  -2.77%  "
$ws.Range("D27").Value = "'153.30"
$ws.Range("E27").Value = "  -1.99%  "
$ws.Range("E28").Value = "  -2.92%  "
$ws.Range("D29").Value = "'6.052"
$ws.Range("E29").Value = "  -4.86%  "
$ws.Range("D30").Value = "'2.009"
$ws.Range("E30").Value = "  -3.74%  "
$ws.Range("D31").Value = "'117.57"
$ws.Range("E31").Value = "  -2.61%  "
$ws.Range("D32").Value = "'0.9667"
$ws.Range("E32").Value = "  -3.58%  "
$ws.Range("E33").Value = "  -2.48%  "
$ws.Range("D34").Value = "'5.365"
$ws.Range("E34").Value = "  -4.72%  "
$ws.Range("D35").Value = "'3.466"
$ws.Range("E35").Value = "  -2.60%  "
$ws.Range("D36").Value = "'1.348"
$ws.Range("E36").Value = "  -3.02%  "
$ws.Range("D37").Value = "'0.06104"
$ws.Range("E37").Value = "  -6.63%  "
$ws.Range("D38").Value = "'0.02201"
$ws.Range("E38").Value = "  -3.53%  "
$ws.Range("D39").Value = "'1.163"
$ws.Range("E39").Value = "  -4.35%  "
$ws.Range("B40").Value = "TheSandbox"
$ws.Range("C40").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D40").Value = "'0.5697"
$ws.Range("E40").Value = "  -3.94%  "
$ws.Range("B41").Value = "FraxShare"
$ws.Range("C41").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D41").Value = "'7.658"
$ws.Range("E41").Value = "  -2.50%  "
$ws.Range("D42").Value = "'10.12"
$ws.Range("E42").Value = "  -5.72%  "
$ws.Range("D43").Value = "'0.1792"
$ws.Range("E43").Value = "  -2.85%  "
$ws.Range("D44").Value = "'2.426"
$ws.Range("E44").Value = "  -2.34%  "
$ws.Range("D45").Value = "'1.247"
$ws.Range("E45").Value = "  +0.42%  "
$ws.Range("D46").Value = "'0.5379"
$ws.Range("E46").Value = "  -3.10%  "
$ws.Range("D47").Value = "'11.72"
$ws.Range("E47").Value = "  -5.59%  "
$ws.Range("D48").Value = "'0.07092"
$ws.Range("E48").Value = "  -6.09%  "
$ws.Range("D49").Value = "'1.901"
$ws.Range("E49").Value = "  -2.97%  "
$ws.Range("D50").Value = "'113.05"
$ws.Range("E50").Value = "  -4.42%  "
$ws.Range("D51").Value = "'2.344"
$ws.Range("E51").Value = "  -3.66%  "
